# New crime data collected — weekly CompStat report refresh
# (Central Park Precinct, week covering 10/23/2023 - 10/29/2023)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: bump the volume/issue number and the reporting week dates ----
# These cells hold multi-run rich text ("Volume 30   Number  42", and the
# "Report Covering the Week  10/16/2023  Through  10/22/2023" banner); use
# Characters() to retarget just the substrings that changed, same as Excel's
# Range.Characters API, so the rest of the text/formatting is left alone.
$ws.Range("A8").Characters(21, 2).Text = "43"
$ws.Range("C9").Characters(27, 10).Text = "10/23/2023"
$ws.Range("C9").Characters(48, 10).Text = "10/29/2023"

# --- Row 15 (Rape) : 2-Year % change recomputed ----------------------------
$ws.Range("N15").Value = -77.777777777777

# --- Row 16 (Robbery) ------------------------------------------------------
$ws.Range("C16").Value = "0"
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = -100
$ws.Range("G16").Value = 5
$ws.Range("H16").Value = -80
$ws.Range("J16").Value = 26
$ws.Range("K16").Value = -30.769230769230
$ws.Range("L16").Value = 28.571428571428
$ws.Range("N16").Value = -89.595375722543

# --- Row 17 (Fel. Assault) --------------------------------------------------
$ws.Range("G17").Value = "0"
$ws.Range("H17").Value = "***.*"
$ws.Range("M17").Value = 60
$ws.Range("N17").Value = -78.378378378378

# --- Row 19 (Gr. Larceny) ---------------------------------------------------
$ws.Range("C19").Value = 1
$ws.Range("D19").Value = "0"
$ws.Range("E19").Value = "***.*"
$ws.Range("F19").Value = 4
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 45
$ws.Range("K19").Value = 95.652173913043
$ws.Range("L19").Value = 104.545454545455
$ws.Range("M19").Value = -23.728813559322
$ws.Range("N19").Value = -71.337579617834

# --- Row 21 (TOTAL) ---------------------------------------------------------
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 8
$ws.Range("H21").Value = -11.111111111111
$ws.Range("I21").Value = 78
$ws.Range("J21").Value = 67
$ws.Range("K21").Value = 16.417910447761
$ws.Range("L21").Value = 52.941176470588
$ws.Range("M21").Value = -17.021276595744
$ws.Range("N21").Value = -81.339712918660

# --- Row 24 (Petit Larceny) --------------------------------------------------
$ws.Range("C24").Value = 1
$ws.Range("D24").Value = 1
$ws.Range("E24").Value = 0
$ws.Range("G24").Value = 4
$ws.Range("H24").Value = -50
$ws.Range("I24").Value = 34
$ws.Range("J24").Value = 28
$ws.Range("K24").Value = 21.428571428571
$ws.Range("L24").Value = 21.428571428571
$ws.Range("M24").Value = -60

# --- Row 25 (Misd. Assault) --------------------------------------------------
$ws.Range("F25").Value = 4
$ws.Range("G25").Value = 1
$ws.Range("H25").Value = 300
$ws.Range("I25").Value = 45
$ws.Range("K25").Value = 66.666666666666
$ws.Range("L25").Value = 104.545454545455
$ws.Range("M25").Value = 150
